# Test.xlsx: "GetUrl" row's value (B2) changes from "Y" to "N".
# This also drops the now-unused "Y" shared-string entry when the
# workbook is re-saved (uniqueCount 6 -> 5), shifting "N" into the
# slot "Y" used to occupy.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = "N"
